$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value (46074 = 2026-02-21) for
# every data row (2..295). Bump it by one day (46075 = 2026-02-22) to match
# the automatic "last changed" stamp refresh recorded in the commit.
for ($r = 2; $r -le 295; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46074) {
        $cell.Value2 = 46075
    }
}
